$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2,  45903, 5588.32624579958, 5071.56761095425, 5760, 6209.32448,  -2.80975645188873),
    @(3,  45904, 5616.87792343629, 5060.12181197823, 2952, 6257.595135, 114.534959314247),
    @(4,  45905, 5616.8779234379,  4345.98004110741, 2952, 6257.595148, 84.7790527362297),
    @(5,  45906, 1140.59185069602, 2176.29459717412, 2952, 2210.680048, 12.2659497699208),
    @(6,  45907, 1021.14313653698, 2124.54061101136, 2952, 2080.094216, 9.64548710309911),
    @(7,  45908, 5899.19381188993, 4915.64258837857, 2952, 6751.2888,   117.32239902036),
    @(8,  45909, 6093.87809712743, 5208.14683736067, 2952, 6751.2888,   121.398230843052),
    @(9,  45910, 6093.87809712743, 5269.19970239235, 2952, 6751.2888,   123.942100219372),
    @(10, 45911, 6093.87809712743, 5258.63445198027, 2952, 6751.2888,   123.501881452202),
    @(11, 45912, 6093.87809712743, 4499.8077894658,  2952, 6751.2888,   91.8841038474324),
    @(12, 45913, 1155.32029950294, 2159.27324039215, 2952, 2198.211777, 10.4235299120501),
    @(13, 45914, 1023.06998026009, 2100.64954231777, 2952, 2053.108927, 7.44535371073668),
    @(14, 45915, 6061.80298841384, 4986.02119019893, 2952, 6762.358872, 113.940711407712),
    @(15, 45916, 6061.80298841384, 5180.34917669728, 2952, 6762.358872, 122.037710845143)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
}
